$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.457.03'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.571.69'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '288.22'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3718'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.82%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.33'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.75%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3317'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.90%  '
$ws.Range("E10").Value = '  -0.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07486'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("E13").Value = '  -1.91%  '
$ws.Range("E14").Value = '  -1.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.903'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.572.21'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("E17").Value = '  +0.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '87.80'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.71%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06741'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.356'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("E22").Value = '  +1.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.447.55'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.381'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.67%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.579'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.29'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.72%  '
$ws.Range("E28").Value = '  -0.84%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.020'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.46'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.747.64'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.053'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.20%  '
$ws.Range("E33").Value = '  -0.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.138'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.783'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08323'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02464'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2271'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.54%  '
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.291'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.351'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6309'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.95%  '
$ws.Range("E43").Value = '  +0.62%  '
$ws.Range("E44").Value = '  -0.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.91'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6156'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.773'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.75'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.210'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07219'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.61%  '
